$d = $word.ActiveDocument

# --- 1. First paragraph: add trailing spaces + red "(This is a change – Version for main branch)" runs ---

# 1a. Append two trailing spaces to the existing (black) run text.
$p1 = $d.Paragraphs.Item(1)
$end1 = $p1.Range.End - 1
$ip = $d.Range($end1, $end1)
$ip.InsertAfter("  ")

# 1b. Insert the first red chunk: "(This is a change \u2013 Ve"
$p1 = $d.Paragraphs.Item(1)
$end1 = $p1.Range.End - 1
$chunk1 = "(This is a change " + [char]0x2013 + " Ve"
$ip = $d.Range($end1, $end1)
$ip.InsertAfter($chunk1)
$r1 = $d.Range($end1, $end1 + $chunk1.Length)
$r1.Font.Color = 255

# 1c. Insert the second red chunk: "rsion for main branch"
$p1 = $d.Paragraphs.Item(1)
$end1 = $p1.Range.End - 1
$chunk2 = "rsion for main branch"
$ip = $d.Range($end1, $end1)
$ip.InsertAfter($chunk2)
$r2 = $d.Range($end1, $end1 + $chunk2.Length)
$r2.Font.Color = 255

# 1d. Insert the third red chunk: ")"
$p1 = $d.Paragraphs.Item(1)
$end1 = $p1.Range.End - 1
$chunk3 = ")"
$ip = $d.Range($end1, $end1)
$ip.InsertAfter($chunk3)
$r3 = $d.Range($end1, $end1 + $chunk3.Length)
$r3.Font.Color = 255

# --- 2. Delete the trailing "ank God almighty, we are free at last." paragraph ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastPara.Range.Delete()
